$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B ("Year of Treatment"); remaining columns C:K shift left to B:J
$ws.Columns.Item(2).Delete()

# Append the ".deja.deja.deja" suffix to the substance-category headers (now in B1:J1)
for ($col = 2; $col -le 10; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $cell.Text + ".deja.deja.deja"
}
